$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily NZ vaccination data rows (8 Nov 2021 - 14 Nov 2021)
$newRows = @(
    @(44508, 5103, 16089),
    @(44509, 5874, 16304),
    @(44510, 6045, 15962),
    @(44511, 6285, 16509),
    @(44512, 6343, 16608),
    @(44513, 7149, 19847),
    @(44514, 4645, 9993)
)

$startRow = 265
$formatSource = $ws.Cells.Item($startRow - 1, 1)
$formatSource.Copy()

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $data[0]
    $dateCell.PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
}

$excel.CutCopyMode = $false

# Update the view to reflect scrolling down to the new bottom rows
$ws.Range("A202").Select()
$excel.ActiveWindow.ScrollRow = 202
$ws.Range("C264").Select()
